$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.144.51'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '2.990.16'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '354.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.46%  '
$ws.Range("E7").Value = '  -0.73%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.628'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.140'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0859'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.94%  '
$ws.Range("D14").Value = '3.460.87'
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.69'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("D16").Value = '2.983.08'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.36%  '
$ws.Range("D18").Value = '52.173.48'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.66'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.42%  '
$ws.Range("D22").Value = '0.0₃0976'
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.181'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  -6.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '36.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.18%  '
$ws.Range("E34").Value = '  -3.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0445'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.21'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.44%  '
$ws.Range("E40").Value = '  -4.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.56%  '
$ws.Range("E45").Value = '  -1.74%  '
$ws.Range("D46").Value = '2.129.02'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.244'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0333'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.932'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.89%  '
